# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "28.264.76", "0.9998") that must stay
# literal text, not be coerced to a number/date by COM value assignment. Force the
# whole data range to Text format first, write the values, then strip the format
# back to Normal so the on-disk style matches the original (no visible NumberFormat).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.264.76"
$ws.Range("E2").Value = "  +3.02%  "

$ws.Range("D3").Value = "1.823.46"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").Value = "340.04"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("D6").Value = "0.9973"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  +3.48%  "

$ws.Range("D8").Value = "0.3502"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").Value = "48.10"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").Value = "1.204"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").Value = "0.07598"
$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").Value = "0.9978"
$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("D13").Value = "22.25"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").Value = "6.554"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").Value = "1.826.37"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "7.218"
$ws.Range("E16").Value = "  +2.23%  "

$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "0.06688"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").Value = "85.60"
$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("D20").Value = "0.9978"
$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "17.90"
$ws.Range("E21").Value = "  +3.32%  "

$ws.Range("D22").Value = "6.606"
$ws.Range("E22").Value = "  +1.68%  "

$ws.Range("D23").Value = "28.237.58"
$ws.Range("E23").Value = "  +2.99%  "

$ws.Range("D24").Value = "12.78"
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("D25").Value = "2.401"
$ws.Range("E25").Value = "  -1.01%  "

$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").Value = "1.492"
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("D28").Value = "21.52"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").Value = "154.85"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").Value = "2.033.58"
$ws.Range("E30").Value = "  +1.62%  "

$ws.Range("D31").Value = "136.14"
$ws.Range("E31").Value = "  +1.75%  "

$ws.Range("D32").Value = "6.237"
$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("D33").Value = "4.040"
$ws.Range("E33").Value = "  -0.55%  "

$ws.Range("D34").Value = "0.08856"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").Value = "13.29"
$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").Value = "5.546"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "0.6986"
$ws.Range("E37").Value = "  +1.45%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02446"
$ws.Range("E38").Value = "  +4.77%  "

$ws.Range("D39").Value = "0.06584"
$ws.Range("E39").Value = "  +3.37%  "

$ws.Range("D40").Value = "1.617"
$ws.Range("E40").Value = "  -4.12%  "

$ws.Range("D41").Value = "0.2231"
$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("D42").Value = "1.272"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").Value = "8.564"
$ws.Range("E43").Value = "  -3.86%  "

$ws.Range("E44").Value = "  +1.16%  "

$ws.Range("D45").Value = "0.6521"
$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").Value = "2.176"
$ws.Range("E47").Value = "  +2.40%  "

$ws.Range("D48").Value = "131.98"
$ws.Range("E48").Value = "  +1.63%  "

$ws.Range("D49").Value = "0.07229"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").Value = "80.74"
$ws.Range("E50").Value = "  +1.77%  "

$ws.Range("B51").Value = "Tezos"
$ws.Range("C51").Value = "https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz"
$ws.Range("D51").Value = "1.248"
$ws.Range("E51").Value = "  +2.81%  "

# Restore plain "Normal" style on column D so no stray number-format survives the edit.
$dRange.Style = "Normal"
